$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.097.63'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '1.636.68'
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("E4").Value = '  -0.20%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '214.39'
$ws.Range("E5").Value = '  +1.68%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.5240'
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  -0.13%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2606'
$ws.Range("E8").Value = '  -0.95%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06305'
$ws.Range("E9").Value = '  +0.07%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '20.65'
$ws.Range("E10").Value = '  -2.71%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07651'
$ws.Range("E11").Value = '  +1.19%  '
$ws.Range("D12").Value = '1.651.20'
$ws.Range("E12").Value = '  -1.18%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.431'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("D14").Value = '1.860.07'
$ws.Range("E14").Value = '  -1.95%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.5505'
$ws.Range("E15").Value = '  -1.29%  '
$ws.Range("D16").Value = '0.0₅8130'
$ws.Range("E16").Value = '  +2.47%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '65.03'
$ws.Range("E17").Value = '  -2.92%  '
$ws.Range("D18").Value = '26.074.92'
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("E19").Value = '  -0.15%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '4.692'
$ws.Range("E20").Value = '  -1.11%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '188.65'
$ws.Range("E21").Value = '  +0.97%  '
$ws.Range("E22").Value = '  -2.08%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '6.154'
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("E24").Value = '  -0.22%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '146.09'
$ws.Range("E25").Value = '  -1.82%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.1217'
$ws.Range("E26").Value = '  -2.51%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '7.419'
$ws.Range("E27").Value = '  -1.48%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '15.84'
$ws.Range("E28").Value = '  -0.82%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.407'
$ws.Range("E29").Value = '  +4.85%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.05934'
$ws.Range("E30").Value = '  -4.98%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.260'
$ws.Range("E31").Value = '  -1.63%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.446'
$ws.Range("E32").Value = '  -1.74%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.408'
$ws.Range("E33").Value = '  -0.15%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.640'
$ws.Range("E34").Value = '  +0.69%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.9898'
$ws.Range("E35").Value = '  -0.59%  '
$ws.Range("E36").Value = '  +1.02%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.396'
$ws.Range("E37").Value = '  -0.52%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.5744'
$ws.Range("E38").Value = '  -4.90%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.01620'
$ws.Range("E39").Value = '  +0.69%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.8571'
$ws.Range("E40").Value = '  -2.22%  '
$ws.Range("E41").Value = '  -0.27%  '
$ws.Range("D42").Value = '1.037.63'
$ws.Range("E42").Value = '  -6.48%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '5.673'
$ws.Range("E43").Value = '  -7.26%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '100.59'
$ws.Range("E44").Value = '  +0.70%  '
$ws.Range("D45").Value = '1.786.47'
$ws.Range("E45").Value = '  -1.88%  '
$ws.Range("E46").Value = '  -2.74%  '
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.003'
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '8.064'
$ws.Range("E49").Value = '  +0.03%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.05169'
$ws.Range("E50").Value = '  -1.26%  '
$ws.Range("E51").Value = '  -0.70%  '
